$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2" = 2.78
    "G2" = 2.98
    "I2" = 2.64
    "M2" = 1.04
    "N2" = 4.9
    "Q2" = 1.65
    "S2" = 2.58
    "T2" = 1.58
    "V2" = 1.64
    "X2" = 26
    "Y2" = 17.5
    "Z2" = 23
    "AB2" = 18.5
    "AC2" = 11.5
    "AF2" = 23
    "AG2" = 16
    "AH2" = 18.5
    "AI2" = 34
    "AK2" = 34
    "AN2" = 22
    "AO2" = 16.5
    "G3" = 12
    "H3" = 1.4
    "L3" = 1.01
    "M3" = 1.05
    "N3" = 3.75
    "O3" = 1.27
    "P3" = 2
    "Q3" = 1.83
    "R3" = 1.38
    "S3" = 3
    "T3" = 2.04
    "U3" = 1.76
    "V3" = 3
    "W3" = 1.1
    "X3" = 970
    "Y3" = 970
    "Z3" = 970
    "AA3" = 970
    "AB3" = 29
    "AC3" = 970
    "AD3" = 970
    "AE3" = 970
    "AF3" = 110
    "AG3" = 38
    "AH3" = 29
    "AI3" = 50
    "AJ3" = 1000
    "AK3" = 1000
    "AL3" = 1000
    "AM3" = 1000
    "AN3" = 1000
    "AO3" = 970
    "G4" = 1.85
    "H4" = 4.8
    "Q4" = 2.1
    "H5" = 1.5
    "I5" = 1.6
    "J5" = 3.95
    "Q5" = 1.99
    "J6" = 4.1
    "Q6" = 1.72
    "H8" = 7
    "O8" = 1.2
    "P8" = 2.64
    "AL8" = 28
    "AO8" = 95
    "F10" = 2.06
    "G10" = 2.18
    "H10" = 4
    "K10" = 3.7
    "P10" = 1.77
    "Q10" = 2.06
    "S10" = 3.9
    "U10" = 1.98
    "W10" = 1.84
    "X10" = 15
    "Y10" = 16
    "AB10" = 10.5
    "AC10" = 9.4
    "AD10" = 20
    "AF10" = 16
    "AG10" = 13
    "AH10" = 24
    "AJ10" = 34
    "AK10" = 30
    "AM10" = 150
    "AN10" = 24
    "F11" = 1.92
    "G11" = 2.34
    "H11" = 3.95
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
